# The contract's closing/signature block has a long run of empty, bold
# 16pt "No Spacing"/"Normal" paragraphs used purely as vertical spacing
# before the signature-line drawing. This trims that run down, removing
# the redundant blank paragraphs while keeping one right after the
# signature sentence and the couple immediately before the signature-line
# drawing intact.

$d = $word.ActiveDocument

# Locate the paragraph that ends the signature sentence
# ("...Tijuana, B.C. a {date}.").
$count = $d.Paragraphs.Count
$sigIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Tijuana, B.C. a*") {
        $sigIdx = $i
        break
    }
}

if ($sigIdx -gt 0) {
    # Walk forward collecting the run of consecutive empty paragraphs
    # (just a paragraph mark, no drawing) that follow the signature line.
    $i = $sigIdx + 1
    $blankIndices = @()
    while ($i -le $d.Paragraphs.Count) {
        $p = $d.Paragraphs.Item($i)
        $isDrawing = $p.Range.WordOpenXML -like "*w:drawing*"
        $isBlank = ($p.Range.Text.Length -eq 1) -and (-not $isDrawing)
        if (-not $isBlank) {
            break
        }
        $blankIndices += $i
        $i = $i + 1
    }

    # Keep the first blank paragraph (right after the signature line) and
    # the last two blank paragraphs (right before the signature-line
    # drawing); delete everything else in between.
    if ($blankIndices.Count -gt 3) {
        $firstKeep = $blankIndices[0]
        $lastDelete = $blankIndices[$blankIndices.Count - 3]

        $startPara = $d.Paragraphs.Item($firstKeep)
        $endPara = $d.Paragraphs.Item($lastDelete)

        $delRange = $d.Range($startPara.Range.End, $endPara.Range.End)
        $delRange.Delete()
    }
}
